$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.102.29"
$ws.Range("E2").Value = "  +0.12%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.546.75"
$ws.Range("E3").Value = "  +2.93%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.78"
$ws.Range("E5").Value = "  +0.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.82"
$ws.Range("E6").Value = "  +3.14%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.586"
$ws.Range("E8").Value = "  -0.39%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.544.68"
$ws.Range("E9").Value = "  +2.84%  "

$ws.Range("E10").Value = "  -0.04%  "

$ws.Range("E11").Value = "  -2.28%  "

$ws.Range("E13").Value = "  +0.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.51"
$ws.Range("E14").Value = "  +3.24%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.001.20"
$ws.Range("E15").Value = "  +2.91%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.031.85"
$ws.Range("E16").Value = "  +0.22%  "

$ws.Range("E17").Value = "  +1.49%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.548.09"
$ws.Range("E18").Value = "  +3.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.40"
$ws.Range("E19").Value = "  +1.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "335.19"
$ws.Range("E20").Value = "  -1.64%  "

$ws.Range("E21").Value = "  +1.60%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.77"
$ws.Range("E22").Value = "  -0.81%  "

$ws.Range("E23").Value = "  -0.25%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.21"
$ws.Range("E24").Value = "  -0.63%  "

$ws.Range("E25").Value = "  +8.98%  "

$ws.Range("E26").Value = "  -1.95%  "

$ws.Range("E27").Value = "  +8.17%  "

$ws.Range("E28").Value = "  +0.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.44"
$ws.Range("E29").Value = "  +4.19%  "

$ws.Range("E30").Value = "  +7.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0817"
$ws.Range("E31").Value = "  +2.25%  "

$ws.Range("E32").Value = "  +0.67%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "175.96"
$ws.Range("E33").Value = "  -0.47%  "

$ws.Range("E34").Value = "  +3.40%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "413.03"
$ws.Range("E35").Value = "  +11.01%  "

$ws.Range("E36").Value = "  +0.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.93"
$ws.Range("E37").Value = "  +0.56%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.39"
$ws.Range("E38").Value = "  -0.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.75"
$ws.Range("E40").Value = "  +3.26%  "

$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.29"
$ws.Range("E42").Value = "  -3.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "152.51"
$ws.Range("E43").Value = "  +2.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.77"
$ws.Range("E44").Value = "  +1.75%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.03"
$ws.Range("E45").Value = "  +2.28%  "

$ws.Range("E46").Value = "  +0.64%  "

$ws.Range("E47").Value = "  +0.41%  "

$ws.Range("E48").Value = "  +1.54%  "

$ws.Range("E49").Value = "  +5.17%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.33"
$ws.Range("E50").Value = "  +1.98%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.77"
$ws.Range("E51").Value = "  +0.64%  "
